# #5: property boat&car done
# Restructure the "汽車" (car) worksheet (sheet3) so it uses the same
# column layout as the "土地" (land) worksheet: name / capacity / owner /
# register_date / register_reason / acquire_value / property_category /
# category / date / legislator_name / legislator_id / source_file / index.

$wb = $excel.ActiveWorkbook
$wsLand = $wb.Worksheets.Item(1)   # 土地 - reference sheet with the full header/column set already present
$wsCar  = $wb.Worksheets.Item(3)   # 汽車

# --- Header row (row 1) -----------------------------------------------
# B1: name, D1: owner, E1: register_date, F1: register_reason,
# G1: acquire_value all already exist verbatim on the land sheet - reuse
# them (and their style) via Copy so no stray shared-string bytes differ.
$wsLand.Range("B1").Copy($wsCar.Range("B1")) | Out-Null
$wsCar.Range("C1").Value = "capacity"
$wsLand.Range("E1").Copy($wsCar.Range("D1")) | Out-Null
$wsLand.Range("F1").Copy($wsCar.Range("E1")) | Out-Null
$wsLand.Range("G1").Copy($wsCar.Range("F1")) | Out-Null
$wsLand.Range("H1").Copy($wsCar.Range("G1")) | Out-Null
$wsLand.Range("I1:O1").Copy($wsCar.Range("H1")) | Out-Null

# --- Data row (row 2) ---------------------------------------------------
# A2 (index), B2 (name/TOYOTA), C2 (capacity/2362), D2 (owner), E2
# (register_date), F2 (register_reason) and G2 (acquire_value) keep their
# existing values - only the trailing columns need to be appended.
$wsLand.Range("I2:N2").Copy($wsCar.Range("H2")) | Out-Null
$wsCar.Range("N2").Value = 34
